$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new data rows (2026/01/14) were inserted right before the existing
# row 619 (the "2026/12/29" block), pushing all subsequent rows down by 2.
$ws.Rows("619:620").Insert()

# Row 619: 2026/01/14, 水, 13:00 slot -> 3, ranking 29
# Leading apostrophe keeps the date-shaped string literal text (matches the
# existing inlineStr/text date cells in column A) instead of Excel's
# auto date-conversion.
$ws.Range("A619").Value = "'2026/01/14"
$ws.Range("B619").Value = "水"
$ws.Range("C619").Value = 3
$ws.Range("D619").Value = 29

# Row 620: 2026/01/14, 水, 16:00 slot -> 6, ranking 31
$ws.Range("A620").Value = "'2026/01/14"
$ws.Range("B620").Value = "水"
$ws.Range("C620").Value = 6
$ws.Range("D620").Value = 31
